$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2230769230769231
$ws.Range("C2").Value = 0.5115384615384615
$ws.Range("J2").Value = 0.01923076923076923
$ws.Range("P2").Value = 0.1653846153846154
$ws.Range("S2").Value = 0.08076923076923077
$ws.Range("B3").Value = 0.01526717557251908
$ws.Range("C3").Value = 0.02290076335877863
$ws.Range("J3").Value = 0.03053435114503817
$ws.Range("P3").Value = 0.7786259541984732
$ws.Range("S3").Value = 0.1526717557251908
$ws.Range("J4").Value = 0.02325581395348837
$ws.Range("P4").Value = 0.7209302325581395
$ws.Range("S4").Value = 0.2558139534883721
$ws.Range("P5").Value = 0.8
$ws.Range("S5").Value = 0.2
$ws.Range("B6").Value = 0.05990783410138249
$ws.Range("D6").Value = 0.02764976958525346
$ws.Range("F6").Value = 0.04147465437788019
$ws.Range("J6").Value = 0.1612903225806452
$ws.Range("O6").Value = 0.004608294930875576
$ws.Range("Q6").Value = 0.1981566820276498
$ws.Range("R6").Value = 0.07834101382488479
$ws.Range("S6").Value = 0.4285714285714285
$ws.Range("B7").Value = 0.08974358974358974
$ws.Range("D7").Value = 0.02564102564102564
$ws.Range("E7").Value = 0.004273504273504274
$ws.Range("F7").Value = 0.03846153846153846
$ws.Range("J7").Value = 0.1538461538461539
$ws.Range("O7").Value = 0.0170940170940171
$ws.Range("Q7").Value = 0.1367521367521368
$ws.Range("R7").Value = 0.1196581196581197
$ws.Range("S7").Value = 0.4145299145299146
$ws.Range("B8").Value = 0.07920792079207921
$ws.Range("D8").Value = 0.01237623762376238
$ws.Range("E8").Value = 0.002475247524752475
$ws.Range("F8").Value = 0.05198019801980198
$ws.Range("J8").Value = 0.1386138613861386
$ws.Range("O8").Value = 0.01485148514851485
$ws.Range("Q8").Value = 0.1559405940594059
$ws.Range("R8").Value = 0.1014851485148515
$ws.Range("S8").Value = 0.4430693069306931
$ws.Range("B9").Value = 0.1097560975609756
$ws.Range("D9").Value = 0.03048780487804878
$ws.Range("E9").Value = 0.006097560975609756
$ws.Range("F9").Value = 0.07317073170731707
$ws.Range("J9").Value = 0.1036585365853658
$ws.Range("Q9").Value = 0.1402439024390244
$ws.Range("R9").Value = 0.07926829268292683
$ws.Range("S9").Value = 0.4573170731707317
$ws.Range("B10").Value = 0.0993676603432701
$ws.Range("D10").Value = 0.01897018970189702
$ws.Range("E10").Value = 0.001806684733514002
$ws.Range("F10").Value = 0.0966576332429991
$ws.Range("J10").Value = 0.1074977416440831
$ws.Range("O10").Value = 0.01987353206865402
$ws.Range("Q10").Value = 0.1725383920505872
$ws.Range("R10").Value = 0.0975609756097561
$ws.Range("S10").Value = 0.3857271906052394
$ws.Range("G11").Value = 0.1424657534246575
$ws.Range("J11").Value = 0.07671232876712329
$ws.Range("K11").Value = 0.2246575342465753
$ws.Range("L11").Value = 0.5315068493150685
$ws.Range("S11").Value = 0.02465753424657534
$ws.Range("G12").Value = 0.7563451776649747
$ws.Range("J12").Value = 0.1928934010152284
$ws.Range("L12").Value = 0.01522842639593909
$ws.Range("S12").Value = 0.03553299492385787
$ws.Range("F13").Value = 0.01886792452830189
$ws.Range("G13").Value = 0.6981132075471698
$ws.Range("J13").Value = 0.1509433962264151
$ws.Range("S13").Value = 0.1320754716981132
$ws.Range("F15").Value = 0.02051282051282051
$ws.Range("H15").Value = 0.1692307692307692
$ws.Range("I15").Value = 0.03076923076923077
$ws.Range("J15").Value = 0.3384615384615385
$ws.Range("K15").Value = 0.06666666666666667
$ws.Range("M15").Value = 0.01538461538461539
$ws.Range("O15").Value = 0.08205128205128205
$ws.Range("S15").Value = 0.2769230769230769
$ws.Range("F16").Value = 0.005847953216374269
$ws.Range("H16").Value = 0.1345029239766082
$ws.Range("I16").Value = 0.08771929824561403
$ws.Range("J16").Value = 0.3684210526315789
$ws.Range("K16").Value = 0.1111111111111111
$ws.Range("M16").Value = 0.02339181286549707
$ws.Range("O16").Value = 0.06432748538011696
$ws.Range("S16").Value = 0.2046783625730994
$ws.Range("F17").Value = 0.0113314447592068
$ws.Range("H17").Value = 0.1558073654390935
$ws.Range("I17").Value = 0.06232294617563739
$ws.Range("J17").Value = 0.424929178470255
$ws.Range("K17").Value = 0.1048158640226629
$ws.Range("M17").Value = 0.0226628895184136
$ws.Range("O17").Value = 0.0424929178470255
$ws.Range("S17").Value = 0.1756373937677054
$ws.Range("F18").Value = 0.009569377990430622
$ws.Range("H18").Value = 0.1435406698564593
$ws.Range("I18").Value = 0.09569377990430622
$ws.Range("J18").Value = 0.2822966507177033
$ws.Range("K18").Value = 0.1244019138755981
$ws.Range("M18").Value = 0.02392344497607655
$ws.Range("N18").Value = 0.004784688995215311
$ws.Range("O18").Value = 0.07177033492822966
$ws.Range("S18").Value = 0.2440191387559809
$ws.Range("F19").Value = 0.01234567901234568
$ws.Range("H19").Value = 0.1902687000726216
$ws.Range("I19").Value = 0.074800290486565
$ws.Range("J19").Value = 0.3231663035584604
$ws.Range("K19").Value = 0.1336238198983297
$ws.Range("M19").Value = 0.02396514161220044
$ws.Range("O19").Value = 0.06172839506172839
$ws.Range("S19").Value = 0.1801016702977487
